$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change L11 formula from =J11*I18 to =K11/I16
$ws.Range("L11").Formula = "=K11/I16"

# Add new row 14: K14 = K15*I16
$ws.Range("K14").Formula = "=K15*I16"

# Add new row 15: K15 = 33/I16
$ws.Range("K15").Formula = "=33/I16"

# Update the selection to N25
$ws.Range("N25").Select()

$wb.Save()
